$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.962.29"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "3.395.80"
$ws.Range("E3").Value = "  -1.44%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.31"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.60"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("D7").Value = "3.396.66"
$ws.Range("E7").Value = "  -1.41%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.476"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("E11").Value = "  -2.48%  "
$ws.Range("E12").Value = "  +1.61%  "
$ws.Range("D13").Value = "3.974.77"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.07"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("E16").Value = "  -1.41%  "
$ws.Range("D17").Value = "3.396.03"
$ws.Range("E17").Value = "  -1.34%  "
$ws.Range("D18").Value = "61.035.92"
$ws.Range("E18").Value = "  -1.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.14"
$ws.Range("E19").Value = "  -3.13%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.85"
$ws.Range("E20").Value = "  -3.49%  "
$ws.Range("E21").Value = "  -4.89%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "383.52"
$ws.Range("E22").Value = "  -4.92%  "
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "74.23"
$ws.Range("E24").Value = "  -0.29%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -4.56%  "
$ws.Range("D27").Value = "3.533.24"
$ws.Range("E27").Value = "  -1.55%  "
$ws.Range("E28").Value = "  -0.77%  "
$ws.Range("E30").Value = "  -3.15%  "
$ws.Range("E31").Value = "  -3.04%  "
$ws.Range("E32").Value = "  -1.21%  "
$ws.Range("E33").Value = "  -2.18%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -1.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.99"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "167.85"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D38").Value = "3.426.36"
$ws.Range("E38").Value = "  -1.35%  "
$ws.Range("E39").Value = "  -2.62%  "
$ws.Range("E40").Value = "  -4.78%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.87"
$ws.Range("E41").Value = "  +2.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0772"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  -1.95%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.67"
$ws.Range("E46").Value = "  -3.68%  "
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").Value = "2.482.15"
$ws.Range("E48").Value = "  -5.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.82"
$ws.Range("E49").Value = "  -2.04%  "
$ws.Range("E50").Value = "  -0.35%  "
$ws.Range("E51").Value = "  +0.99%  "

$wb.Save()